$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1335
$ws.Range("I32").Value = 337.5
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 337.5
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -11.5
$ws.Range("N32").Value = -2652

$ws.Range("H129").Value = 737.1
$ws.Range("I129").Value = 381.42856
$ws.Range("J129").Value = 928.61536
$ws.Range("K129").Value = 1144.28568
$ws.Range("L129").Value = 2785.84608
$ws.Range("M129").Value = 3855.71432
$ws.Range("N129").Value = -12785.84608

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1203.9445
$ws.Range("I2").Value = 1224.7273
$ws.Range("K2").Value = 1224.7273
$ws.Range("M2").Value = -1111.7273

$ws.Range("H32").Value = 5710.2627
$ws.Range("I32").Value = 4662.716
$ws.Range("J32").Value = 14090.637
$ws.Range("K32").Value = 4662.716
$ws.Range("L32").Value = 14090.637
$ws.Range("M32").Value = -4375.716
$ws.Range("N32").Value = -14664.637

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H102").Value = 10417734
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws.Range("H116").Value = 1203.9445
$ws.Range("I116").Value = 1224.7273
$ws.Range("K116").Value = 1224.7273
$ws.Range("M116").Value = 1069.2727

$ws.Range("H125").Value = 35657.5
$ws.Range("J125").Value = 35657.5
$ws.Range("L125").Value = 35657.5
$ws.Range("N125").Value = -45497.5

$ws.Range("H132").Value = 2179.0833
$ws.Range("I132").Value = 1843.8276
$ws.Range("J132").Value = 2690.7896
$ws.Range("K132").Value = 5531.4828
$ws.Range("L132").Value = 8072.3688
$ws.Range("M132").Value = -3001.4828
$ws.Range("N132").Value = -13132.3688

$ws.Range("H134").Value = 31360
$ws.Range("J134").Value = 31360
$ws.Range("L134").Value = 31360
$ws.Range("N134").Value = -41500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1203.9445
$ws.Range("I3").Value = 1224.7273
$ws.Range("K3").Value = 1224.7273
$ws.Range("M3").Value = -1110.7273

$ws.Range("H107").Value = 1436.6666
$ws.Range("I107").Value = 1079.8334
$ws.Range("J107").Value = 2150.3333
$ws.Range("K107").Value = 1079.8334
$ws.Range("L107").Value = 2150.3333
$ws.Range("M107").Value = 840.1666
$ws.Range("N107").Value = -5990.3333

$ws.Range("H134").Value = 880.4643
$ws.Range("I134").Value = 755.9231
$ws.Range("K134").Value = 2267.7693
$ws.Range("M134").Value = 267.2307000000001

$ws.Range("H139").Value = 33002.5
$ws.Range("J139").Value = 33002.5
$ws.Range("L139").Value = 33002.5
$ws.Range("N139").Value = -43282.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 47272.25
$ws.Range("J20").Value = 47272.25
$ws.Range("L20").Value = 47272.25
$ws.Range("N20").Value = -47744.25

$ws.Range("H30").Value = 47272.25
$ws.Range("J30").Value = 47272.25
$ws.Range("L30").Value = 47272.25
$ws.Range("N30").Value = -47454.25

$ws.Range("H86").Value = 3367743
$ws.Range("I86").Value = 4466523
$ws.Range("K86").Value = 4466523
$ws.Range("M86").Value = -4465400

$ws.Range("H89").Value = 3367743
$ws.Range("I89").Value = 4466523
$ws.Range("K89").Value = 22332615
$ws.Range("M89").Value = -22326999

$ws.Range("H122").Value = 846.4737
$ws.Range("I122").Value = 782.3889
$ws.Range("K122").Value = 2347.1667
$ws.Range("M122").Value = 102.8332999999998

$ws.Range("H128").Value = 47272.25
$ws.Range("J128").Value = 47272.25
$ws.Range("L128").Value = 47272.25
$ws.Range("N128").Value = -57232.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10953
$ws.Range("I3").Value = 6529
$ws.Range("K3").Value = 19587
$ws.Range("M3").Value = -19475

$ws.Range("H5").Value = 1014.5455
$ws.Range("I5").Value = 1128.9584
$ws.Range("J5").Value = 709.44446
$ws.Range("K5").Value = 3386.8752
$ws.Range("L5").Value = 2128.33338
$ws.Range("M5").Value = -3274.8752
$ws.Range("N5").Value = -2352.33338

$ws.Range("H131").Value = 21277628
$ws.Range("I131").Value = 142857700
$ws.Range("J131").Value = 1115.85
$ws.Range("K131").Value = 428573100
$ws.Range("L131").Value = 3347.55
$ws.Range("M131").Value = -428568060
$ws.Range("N131").Value = -13427.55

$ws.Range("H135").Value = 1014.5455
$ws.Range("I135").Value = 1128.9584
$ws.Range("J135").Value = 709.44446
$ws.Range("K135").Value = 10160.6256
$ws.Range("L135").Value = 6385.00014
$ws.Range("M135").Value = -7625.625599999999
$ws.Range("N135").Value = -11455.00014

$ws.Range("H140").Value = 2821.353
$ws.Range("I140").Value = 1832.6111
$ws.Range("J140").Value = 3360.6667
$ws.Range("K140").Value = 5497.8333
$ws.Range("L140").Value = 10082.0001
$ws.Range("M140").Value = -317.8333000000002
$ws.Range("N140").Value = -20442.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 112502490
$ws.Range("I70").Value = 83336650
$ws.Range("J70").Value = 200000000
$ws.Range("K70").Value = 83336650
$ws.Range("L70").Value = 200000000
$ws.Range("M70").Value = -83336380
$ws.Range("N70").Value = -200000540

$ws.Range("H73").Value = 112502490
$ws.Range("I73").Value = 83336650
$ws.Range("J73").Value = 200000000
$ws.Range("K73").Value = 83336650
$ws.Range("L73").Value = 200000000
$ws.Range("M73").Value = -83335714
$ws.Range("N73").Value = -200001872

$ws.Range("H97").Value = 681.8182
$ws.Range("I97").Value = 686.4286
$ws.Range("J97").Value = 673.75
$ws.Range("K97").Value = 686.4286
$ws.Range("L97").Value = 673.75
$ws.Range("M97").Value = -190.4286
$ws.Range("N97").Value = -1665.75

$ws.Range("H113").Value = 1234.2273
$ws.Range("I113").Value = 1116.8334
$ws.Range("J113").Value = 1762.5
$ws.Range("K113").Value = 1116.8334
$ws.Range("L113").Value = 1762.5
$ws.Range("M113").Value = 1053.1666
$ws.Range("N113").Value = -6102.5

$ws.Range("H126").Value = 1803
$ws.Range("I126").Value = 1442.091
$ws.Range("K126").Value = 4326.272999999999
$ws.Range("M126").Value = -1856.272999999999

$ws.Range("H132").Value = 2882.4138
$ws.Range("I132").Value = 2643.2
$ws.Range("J132").Value = 3414
$ws.Range("K132").Value = 7929.599999999999
$ws.Range("L132").Value = 10242
$ws.Range("M132").Value = -5399.599999999999
$ws.Range("N132").Value = -15302

$ws.Range("H135").Value = 34173.57
$ws.Range("J135").Value = 33587.406
$ws.Range("L135").Value = 33587.406
$ws.Range("N135").Value = -43727.406

$ws.Range("H141").Value = 38704.145
$ws.Range("J141").Value = 38704.145
$ws.Range("L141").Value = 38704.145
$ws.Range("N141").Value = -49064.145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4994.5
$ws.Range("I40").Value = 2756.3
$ws.Range("K40").Value = 2756.3
$ws.Range("M40").Value = -2620.3

$ws.Range("H61").Value = 1039.1818
$ws.Range("I61").Value = 980.6667
$ws.Range("J61").Value = 1302.5
$ws.Range("K61").Value = 980.6667
$ws.Range("L61").Value = 1302.5
$ws.Range("M61").Value = -778.6667
$ws.Range("N61").Value = -1706.5

$ws.Range("H113").Value = 1039.1818
$ws.Range("I113").Value = 980.6667
$ws.Range("J113").Value = 1302.5
$ws.Range("K113").Value = 980.6667
$ws.Range("L113").Value = 1302.5
$ws.Range("M113").Value = 1189.3333
$ws.Range("N113").Value = -5642.5

$ws.Range("H132").Value = 2970.2222
$ws.Range("I132").Value = 2743.2
$ws.Range("K132").Value = 8229.599999999999
$ws.Range("M132").Value = -5699.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 449.23077
$ws.Range("I107").Value = 325.75
$ws.Range("J107").Value = 646.8
$ws.Range("K107").Value = 977.25
$ws.Range("L107").Value = 1940.4
$ws.Range("M107").Value = 942.75
$ws.Range("N107").Value = -5780.4

$ws.Range("H122").Value = 13159453
$ws.Range("I122").Value = 13890422
$ws.Range("K122").Value = 41671266
$ws.Range("M122").Value = -41668816

$ws.Range("H126").Value = 62500424
$ws.Range("I126").Value = 62500424
$ws.Range("K126").Value = 187501272
$ws.Range("M126").Value = -187498802

$ws.Range("H131").Value = 35405
$ws.Range("J131").Value = 35405
$ws.Range("L131").Value = 35405
$ws.Range("N131").Value = -45485
